# Adds summary rows (average/worst of SW and SC ratios) below the
# per-instance data table on Sheet1, matching the author's manual
# addition of:
#   J12  = AVERAGE(J2:J11)
#   A14  = "Average of SW(S*)/SW(OPT)"   B14 = AVERAGE(N2:N11)
#   A15  = "Average of SC(S*)/SC(OPT)"   B15 = AVERAGE(Z2:Z11)
#   A16  = "Worst of SW(S*)/SW(OPT)"     B16 = MIN(N2:N11)
#   A17  = "Worst of SC(S*)/SC(OPT)"     B17 = MAX(Z2:Z11)
# with the B14:B17 labels bold/12pt/vertically centred, plus the
# row-height, selection and page-setup touch-ups that come along with
# it when the sheet was last saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Quick overall average (row 12) of the k column.
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# Summary block.
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Bold / 12pt / vertically-centred look for the summary values. Build
# it once on B14 then fan it out with a formats-only paste so the style
# table collapses to a single new font + single new cell style instead
# of accumulating one extra style per property touched.
$ws.Range("B14").Font.Bold = $true
$ws.Range("B14").Font.Size = 12
$ws.Range("B14").VerticalAlignment = -4108
$ws.Range("B14").Copy()
$ws.Range("B15:B17").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Match the slightly taller rows Excel used for the bigger font.
$ws.Rows.Item(14).RowHeight = 15.6
$ws.Rows.Item(15).RowHeight = 15.6
$ws.Rows.Item(16).RowHeight = 15.6
$ws.Rows.Item(17).RowHeight = 15.6

# Leave the selection where the author left it after typing the block.
$ws.Range("A14:B17").Select()

# Page setup touch-up that came with the resave.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
